$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 417
$wsExhibit.Range("F3").Value = 563
$wsExhibit.Range("F5").Value = 608
$wsExhibit.Range("F6").Value = 786
$wsExhibit.Range("F8").Value = 550
$wsExhibit.Range("F10").Value = 1115
$wsExhibit.Range("F11").Value = 567
$wsExhibit.Range("F12").Value = 329
$wsExhibit.Range("F14").Value = 139
$wsExhibit.Range("F15").Value = 288
$wsExhibit.Range("F18").Value = 523
$wsExhibit.Range("F19").Value = 511
$wsExhibit.Range("F21").Value = 457

# Sheet "演出" (sheet2)
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 65
$wsShow.Range("F9").Value = 187
$wsShow.Range("F12").Value = 18

# Sheet "全部类型" (sheet4)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 65
$wsAll.Range("F4").Value = 417
$wsAll.Range("F7").Value = 563
$wsAll.Range("F9").Value = 608
$wsAll.Range("F10").Value = 786
$wsAll.Range("F12").Value = 550
$wsAll.Range("F14").Value = 1115
$wsAll.Range("F15").Value = 567
$wsAll.Range("F18").Value = 329
$wsAll.Range("F21").Value = 139
$wsAll.Range("F23").Value = 288
$wsAll.Range("F26").Value = 187
$wsAll.Range("F28").Value = 523
$wsAll.Range("F30").Value = 18
$wsAll.Range("F31").Value = 511
$wsAll.Range("F33").Value = 457
